$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "30.313 (16.77)" "30.313 (13.86)"
Replace-Text "25.280 (13.99)" "25.280 (11.56)"
Replace-Text "16.853 (9.32)" "16.853 (7.71)"
Replace-Text "40.324 (22.31)" "40.324 (18.44)"
Replace-Text "16.463 (9.11)" "16.463 (7.53)"
Replace-Text "10.145 (5.61)" "10.145 (4.64)"
Replace-Text "17.225 (9.53)" "17.225 (7.88)"
Replace-Text "5.570 (3.08)" "5.570 (2.55)"
Replace-Text "13.835 (7.65)" "13.835 (6.33)"

Replace-Text "Materiais" "Direito"
Replace-Text "4.751 (2.63)" "42.682 (19.52)"
Replace-Text "3.606 (75.9)" "31.722 (74.32)"
Replace-Text "1.145 (24.1)" "10.960 (25.68)"
Replace-Text "2.716 (57.17)" "24.320 (56.98)"
Replace-Text "2.035 (42.83)" "18.362 (43.02)"
Replace-Text "613 (12.9)" "5.409 (12.67)"
Replace-Text "532 (11.2)" "5.551 (13.01)"
Replace-Text "1.422 (29.93)" "12.953 (30.35)"
Replace-Text "2.184 (45.97)" "18.769 (43.97)"
